$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for rows 2-25 (bus indices 0-23), columns B-F and I-M
# (slack bus voltage setpoint changed from 1.05 to 1.02 -> all dependent bus voltages recomputed)
$data = @{
    2 = @(1.02, 1.085301709118959, 1.085124543965622, 1.087265549381269, 1.096873305803911, 1.059245209726623, 1.090157246685576, 1.087785785571798, 1.089921258391642, 1.099504473980933)
    3 = @(1.02, 1.08701399729553, 1.086493387049638, 1.088776999058817, 1.098384363779024, 1.059759922991254, 1.091529202272312, 1.088972439691653, 1.091250581336528, 1.100835199105746)
    4 = @(1.02, 1.088119578485727, 1.087376886973792, 1.089752568116499, 1.099359950300614, 1.060090394316215, 1.092414179809803, 1.089737510452409, 1.092107789751352, 1.101693581356675)
    5 = @(1.02, 1.088583806880934, 1.087747785041195, 1.090162124136735, 1.099769577973882, 1.060228710566108, 1.092785571623838, 1.090058490920428, 1.092467463566308, 1.102053811967879)
    6 = @(1.02, 1.088661720521898, 1.087810029923799, 1.090230857096993, 1.099838326720331, 1.060251898600485, 1.0928478919635, 1.090112346781096, 1.09252781382869, 1.102114259417345)
    7 = @(1.02, 1.088125783703029, 1.087381844983109, 1.08975804286134, 1.099365425750921, 1.060092244910822, 1.092419144916046, 1.089741801971794, 1.092112598453544, 1.101698397250679)
    8 = @(1.02, 1.085880886192941, 1.08558761889134, 1.087776862679384, 1.097384430298241, 1.059419696870241, 1.090621485142614, 1.088187400524059, 1.090371126791931, 1.099954760770666)
    9 = @(1.02, 1.081906259846066, 1.082408441323232, 1.08426659148753, 1.093876584021221, 1.058214604362876, 1.087432105509631, 1.085426707197246, 1.087279359456254, 1.096861242041419)
    10 = @(1.02, 1.079243043975717, 1.080276619470768, 1.081912826054108, 1.091525880624595, 1.057397499911419, 1.085290630642364, 1.083571140111883, 1.085202023534982, 1.094784156202508)
    11 = @(1.02, 1.078086468508088, 1.079350450965541, 1.080890243180368, 1.090504973939937, 1.057040370504787, 1.08435959254955, 1.082763951871856, 1.084298537337753, 1.093881121233626)
    12 = @(1.02, 1.077656340680966, 1.079005957199657, 1.080509888140073, 1.090125294777881, 1.056907212825952, 1.084013185225237, 1.082463557190281, 1.083962330947784, 1.093545134544137)
    13 = @(1.02, 1.077748628533043, 1.079079873896297, 1.08059149943573, 1.090206758674655, 1.056935798491562, 1.08408751714652, 1.082528018793968, 1.084034476228866, 1.093617230332157)
    14 = @(1.02, 1.078050924791968, 1.079321984730382, 1.080858813645339, 1.090473599182888, 1.057029373961024, 1.084330970283131, 1.082739132824818, 1.08427075895776, 1.093853359948128)
    15 = @(1.02, 1.0782371095635, 1.079471094111213, 1.081023445149018, 1.090637945901019, 1.057086961958141, 1.084480892782466, 1.08286913141919, 1.084416259153383, 1.093998772725929)
    16 = @(1.02, 1.079319729427141, 1.080338020448661, 1.081980619013867, 1.091593569776106, 1.057421131047643, 1.085352340168404, 1.083624631335349, 1.085261899890331, 1.094844009885398)
    17 = @(1.02, 1.079997910940931, 1.080880988954003, 1.082580112711142, 1.09219218544192, 1.057629854389641, 1.08589795913022, 1.084097534134986, 1.085791272000037, 1.095373221059258)
    18 = @(1.02, 1.080393157160328, 1.081197397330002, 1.082929461365663, 1.092541055712002, 1.057751279437117, 1.086215846997694, 1.084373012587569, 1.086099662053731, 1.095681550292049)
    19 = @(1.02, 1.080527871151404, 1.08130523455359, 1.083048525357452, 1.092659962293468, 1.057792628169923, 1.086324177448964, 1.084466883266056, 1.086204750380478, 1.095786623469056)
    20 = @(1.02, 1.079925182271966, 1.080822764268153, 1.082515826514642, 1.092127989953696, 1.0576074934627, 1.085839456959564, 1.084046833222033, 1.08573451516535, 1.095316478068796)
    21 = @(1.02, 1.077961920674695, 1.079250702289126, 1.080780110797525, 1.090395034352761, 1.057001832272235, 1.084259295468838, 1.082676980800104, 1.084201196561876, 1.093783841197532)
    22 = @(1.02, 1.076724499485266, 1.078259538968714, 1.079685769243856, 1.089302736697671, 1.056618110794992, 1.083262432884125, 1.081812402541452, 1.083233594694755, 1.092816969647702)
    23 = @(1.02, 1.077380773335791, 1.078785237507553, 1.08026619194792, 1.089882046436311, 1.056821807263201, 1.083791210611446, 1.082271048059054, 1.083746878860736, 1.093329837871186)
    24 = @(1.02, 1.079958046232838, 1.080849074411967, 1.082544875686569, 1.092156998034188, 1.057617598391061, 1.085865892692159, 1.084069743887604, 1.085760162321553, 1.095342118867299)
    25 = @(1.02, 1.082936109876331, 1.08323247168638, 1.085176424139053, 1.094785537499505, 1.058528545480589, 1.088259271521976, 1.086143033851411, 1.088081455635591, 1.097663543357359)
}

$cols = @("B","C","D","E","F","I","J","K","L","M")

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}
